# The sheet holds monthly index values for 2014-2017, one 12-row block per
# year (columns A=month label, B:G = series values), stacked in calendar
# order (Jan..Dec) starting at row 2.
#
# The edit re-orders the rows *within each year's 12-row block* so that
# October, November and December come first, followed by January..September,
# i.e. old block order [Jan..Dec] -> new block order [Oct,Nov,Dec,Jan..Sep].
# This is the same as rotating each 12-row block by 3 (each row's new
# position = (old position - 9) mod 12, 0-based).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow  = 2   # first data row (row 1 is the header)
$blockSize = 12  # months per year
$numBlocks = 4   # 2014, 2015, 2016, 2017
$numCols   = 7   # A:G

for ($b = 0; $b -lt $numBlocks; $b++) {
    $blockStart = $startRow + ($b * $blockSize)
    $blockEnd   = $blockStart + $blockSize - 1

    $rng = $ws.Range($ws.Cells.Item($blockStart, 1), $ws.Cells.Item($blockEnd, $numCols))
    $v = $rng.Value()

    $newArr = New-Object 'object[,]' $blockSize, $numCols
    for ($i = 1; $i -le $blockSize; $i++) {
        $oldIndex0 = $i - 1
        $newIndex0 = (($oldIndex0 - 9) % 12 + 12) % 12
        for ($j = 1; $j -le $numCols; $j++) {
            $newArr[$newIndex0, $j - 1] = $v[$i, $j]
        }
    }

    $rng.Value = $newArr
}
